# Fruta / hortaliza, semanal
# Update the weekly price rows (2-17) on the single worksheet: the data for
# "Hortaliza, Femacal de La Calera - Pepino dulce" was re-sorted/refreshed,
# so for each row we rewrite Fecha (D), Calidad (I), Volumen (J),
# Precio minimo/maximo/promedio (K/L/M), Unidad de comercializacion (N),
# Precio $/Kg (P) and Kg o Unidades (Q) to their new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2;  D = 44991; I = "Primera"; J = 75; K = 13000; L = 13000; M = 13000; N = "`$/bandeja 18 kilos"; P = 722;  Q = 18 },
    @{ Row = 3;  D = 44991; I = "Segunda"; J = 56; K = 9000;  L = 9000;  M = 9000;  N = "`$/bandeja 18 kilos"; P = 500;  Q = 18 },
    @{ Row = 4;  D = 44756; I = "Primera"; J = 65; K = 14000; L = 14000; M = 14000; N = "`$/caja 15 kilos";    P = 933;  Q = 15 },
    @{ Row = 5;  D = 44756; I = "Segunda"; J = 68; K = 12000; L = 12000; M = 12000; N = "`$/caja 15 kilos";    P = 800;  Q = 15 },
    @{ Row = 6;  D = 44536; I = "Primera"; J = 87; K = 22000; L = 22000; M = 22000; N = "`$/bandeja 18 kilos"; P = 1222; Q = 18 },
    @{ Row = 7;  D = 44536; I = "Segunda"; J = 80; K = 20000; L = 20000; M = 20000; N = "`$/bandeja 18 kilos"; P = 1111; Q = 18 },
    @{ Row = 8;  D = 44235; I = "Primera"; J = 80; K = 14000; L = 14000; M = 14000; N = "`$/bandeja 18 kilos"; P = 778;  Q = 18 },
    @{ Row = 9;  D = 44235; I = "Segunda"; J = 70; K = 12000; L = 12000; M = 12000; N = "`$/bandeja 18 kilos"; P = 667;  Q = 18 },
    @{ Row = 10; D = 44235; I = "Tercera"; J = 60; K = 10000; L = 10000; M = 10000; N = "`$/bandeja 18 kilos"; P = 556;  Q = 18 },
    @{ Row = 11; D = 44238; I = "Primera"; J = 90; K = 13000; L = 13000; M = 13000; N = "`$/bandeja 18 kilos"; P = 722;  Q = 18 },
    @{ Row = 12; D = 44238; I = "Segunda"; J = 80; K = 11000; L = 11000; M = 11000; N = "`$/bandeja 18 kilos"; P = 611;  Q = 18 },
    @{ Row = 13; D = 44424; I = "Primera"; J = 75; K = 18000; L = 18000; M = 18000; N = "`$/caja 15 kilos";    P = 1200; Q = 15 },
    @{ Row = 14; D = 44424; I = "Segunda"; J = 50; K = 12000; L = 12000; M = 12000; N = "`$/caja 15 kilos";    P = 800;  Q = 15 },
    @{ Row = 15; D = 44242; I = "Primera"; J = 60; K = 13000; L = 13000; M = 13000; N = "`$/bandeja 18 kilos"; P = 722;  Q = 18 },
    @{ Row = 16; D = 44242; I = "Segunda"; J = 50; K = 10000; L = 10000; M = 10000; N = "`$/bandeja 18 kilos"; P = 556;  Q = 18 },
    @{ Row = 17; D = 44992; I = "Primera"; J = 56; K = 13000; L = 13000; M = 13000; N = "`$/bandeja 18 kilos"; P = 722;  Q = 18 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("D$n").Value = $r.D
    $ws.Range("I$n").Value = $r.I
    $ws.Range("J$n").Value = $r.J
    $ws.Range("K$n").Value = $r.K
    $ws.Range("L$n").Value = $r.L
    $ws.Range("M$n").Value = $r.M
    $ws.Range("N$n").Value = $r.N
    $ws.Range("P$n").Value = $r.P
    $ws.Range("Q$n").Value = $r.Q
}
